$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 to the name "hossam" and B2 to the new ssn value
$ws.Range("A2").Value = "hossam"
$ws.Range("B2").Value = 123456789053

# Remove the old row 3 entirely (A3/B3)
$ws.Rows("3:3").Delete()

# Widen column B slightly (Excel stores width in 1/256-character units and
# rounds up, so 11.1 round-trips through the OOXML writer as exactly 12)
$ws.Columns("B:B").ColumnWidth = 11.1

# Update the active selection shown in the saved view
$ws.Range("B2").Select()
